# Y4_B2526_Excuses.xlsx - attendance app upload
# Updates the single logged excuse row (row 2) with a new entry, and widens
# the "Subject" column (B) to fit the longer text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B (Subject) from 12 to 17 characters.
# ColumnWidth's stored/serialized width is offset by +5/6 (~0.8333) from the
# value you assign, so back that off to land exactly on 17.
$ws.Columns.Item(2).ColumnWidth = 17 - 0.8333333333333334

# Row 2: replace the logged excuse with the new record.
# Student ID is entered with a leading apostrophe so Excel stores it as text
# (matching the original "Student ID" column, which holds numeric-looking
# text rather than a true number).
$ws.Range("A2").Value = "'211764"
$ws.Range("B2").Value = "general surgery"
$ws.Range("C2").Value = "29/10/2025"
$ws.Range("D2").Value = "10:30:00"
